# Update workbook for "new runs" data:
#  - Platform Coverage sheet: MDA age-band 15-50 is extended to 15-65
#    (G4: 50 -> 65), which makes the old 50-65 MDA row (row 5) redundant,
#    so that row is removed entirely (all rows below shift up by one).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Platform Coverage")

# Extend the "MDA" (age 15-?) row's max age from 50 to 65.
$ws.Cells.Item(4, 7).Value = 65

# Remove the now-duplicate "MDA" 50-65 row (row 5); everything below
# shifts up by one row.
$ws.Rows(5).EntireRow.Delete()

# Restore the cursor/selection position recorded in the saved file.
$ws.Range("I9").Select()
